$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.273.68'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.60%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.878.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.83%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.12%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -3.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06590'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.52%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.893.97'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.29%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.72'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07288'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.26%  '

$ws.Range('E13').Value = '  +1.47%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6552'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.99%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.230.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9995'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007694'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.31%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.118.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.314'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.0000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.07%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '195.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.24%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.119'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.70%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.292'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.45%  '

$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('E27').Value = '  -4.57%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.914'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.31%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.442'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.75%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.270'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09130'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.067'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.29%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05105'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.43%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7182'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.31%  '

$ws.Range('E35').Value = '  -1.65%  '

$ws.Range('E36').Value = '  +0.91%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01796'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.98%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.640'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9189'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.89%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.13'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4286'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.70%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.808'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9987'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.22%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.401'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.10%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1319'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.09%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.145'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.95%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.98'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.64%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05746'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3824'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.10%  '
